$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: find the paragraph whose text starts with a given prefix.
# ---------------------------------------------------------------------------
function Get-ParagraphByPrefix($doc, $prefix) {
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $para = $doc.Paragraphs($i)
        if ($para.Range.Text.StartsWith($prefix)) {
            return $para
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# 1) Right-align the "Currently, I am interested in the intersection..."
#    paragraph (adds <w:jc w:val="right"/> to its pPr).
# ---------------------------------------------------------------------------
$pAlign = Get-ParagraphByPrefix $d "Currently, I am interested in the intersection"
$pAlign.Range.ParagraphFormat.Alignment = 2

# ---------------------------------------------------------------------------
# 2) Split the "I am a Junior Mechanical Engineering..." paragraph's single
#    run into several runs, interspersed with proofing-error markers
#    (<w:proofErr .../>), without altering the visible text.
# ---------------------------------------------------------------------------
$rPr = '<w:rPr><w:rFonts w:ascii="Helvetica Neue" w:eastAsia="Times New Roman" w:hAnsi="Helvetica Neue" w:cs="Times New Roman"/><w:color w:val="F2F2F2"/><w:sz w:val="30"/><w:szCs w:val="30"/><w:shd w:val="clear" w:color="auto" w:fill="C9B783"/></w:rPr>'

$para5Runs = '<w:r w:rsidRPr="009368AC">' + $rPr + '<w:t xml:space="preserve">I am a Junior Mechanical Engineering student at the New York University. I am passionate about designing, prototyping, and implementing systems that solves a </w:t></w:r>' + `
    '<w:proofErr w:type="gramStart"/>' + `
    '<w:r>' + $rPr + '<w:t>particular problem</w:t></w:r>' + `
    '<w:proofErr w:type="gramEnd"/>' + `
    '<w:r>' + $rPr + '<w:t xml:space="preserve">; and through participating in projects and competitions, I have honed my skills in Arduino, </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r>' + $rPr + '<w:t>Solidworks</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r>' + $rPr + '<w:t xml:space="preserve">, Python, Prototyping and </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r>' + $rPr + '<w:t>Papercrafting</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r>' + $rPr + '<w:t>.</w:t></w:r>'

$para5Xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $para5Runs + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$pJunior = Get-ParagraphByPrefix $d "I am a Junior Mechanical Engineering student"
$fullRange = $pJunior.Range
$textRange = $d.Range($fullRange.Start, $fullRange.End - 1)
$textRange.InsertXML($para5Xml)

# ---------------------------------------------------------------------------
# 3) Split the "Currently, I am working on developing..." paragraph's single
#    run into several runs, interspersed with proofing-error markers,
#    without altering the visible text.
# ---------------------------------------------------------------------------
$para7Runs = '<w:r w:rsidRPr="009368AC">' + $rPr + '<w:t xml:space="preserve">Currently, I am working on developing an Autonomous Surface Vehicle that could gather water quality data autonomously in the Gowanus Canal. Several of my work had gained international recognition, such as the </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r>' + $rPr + '<w:t>Wadi</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r>' + $rPr + '<w:t xml:space="preserve"> Drone. </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r>' + $rPr + '<w:t>Wadi</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r>' + $rPr + '<w:t xml:space="preserve"> Drone is the winning drone of the 2015 Drones </w:t></w:r>' + `
    '<w:proofErr w:type="gramStart"/>' + `
    '<w:r>' + $rPr + '<w:t>For</w:t></w:r>' + `
    '<w:proofErr w:type="gramEnd"/>' + `
    '<w:r>' + $rPr + '<w:t xml:space="preserve"> Good Competition that aids conservation efforts in the UAE.</w:t></w:r>'

$para7Xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $para7Runs + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$pWadi = Get-ParagraphByPrefix $d "Currently, I am working on developing"
$fullRange2 = $pWadi.Range
$textRange2 = $d.Range($fullRange2.Start, $fullRange2.End - 1)
$textRange2.InsertXML($para7Xml)

Write-Output "Applied cornstalk-section edits."
